$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (G1: "Duo/Team" -> "Duo"; M1: trim " Screenshot of Payment " -> "Screenshot of Payment")
$ws.Range("G1").Value = "Duo"
$ws.Range("M1").Value = "Screenshot of Payment"

# Update sheet view: zoom + new selection
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("M2").Select()
